# Bulk upload excel update
# The "brand_id" column (column D) is removed entirely from the product
# bulk-upload template, shifting every later column one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "brand_id" column (D). EntireColumn.Delete shifts all
# columns to the right of D left by one, which matches the diff: the sheet
# dimension shrinks from A1:Z2 to A1:Y2, the brand_id shared string is
# dropped, and every header/value previously in columns E:Z now lives one
# column earlier (D:Y).
$ws.Range("D1").EntireColumn.Delete()

# The old "meta_description" column (was X, now the new column 24 after the
# shift) had become a leftover very wide column (previously width 126+,
# inherited from the old "photos" column Y that sat in that position) -
# narrow it back down to a normal width. The trailing column (now Y /
# "photos") goes back to the sheet's default width (no explicit override).
$ws.Columns.Item(24).ColumnWidth = 26.3

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("B9").Select()
